$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C23").Value = [double]"1.910044682548787e-282"
$ws.Range("C24").Value = [double]"2.579144995786488e-236"
$ws.Range("C25").Value = [double]"2.490596052932349e-203"
$ws.Range("C26").Value = [double]"1.478079197623186e-178"
$ws.Range("C27").Value = [double]"2.966295466886597e-159"
$ws.Range("C28").Value = [double]"8.861400161105788e-144"
$ws.Range("C29").Value = [double]"4.372921143488862e-131"
$ws.Range("C30").Value = [double]"1.770189961932425e-120"
$ws.Range("C31").Value = [double]"1.779584499825226e-111"
$ws.Range("C32").Value = [double]"9.787707032228842e-104"
$ws.Range("C33").Value = [double]"5.251359530553954e-97"
$ws.Range("C34").Value = [double]"4.239381583284432e-91"
$ws.Range("C35").Value = [double]"3.291821495500461e-81"
$ws.Range("C36").Value = [double]"3.002684247432046e-73"
$ws.Range("C37").Value = [double]"8.801468397002664e-59"
$ws.Range("C38").Value = [double]"5.318653252209473e-49"
$ws.Range("C39").Value = [double]"6.634539932282144e-42"
$ws.Range("C40").Value = [double]"1.839717041377047e-36"
$ws.Range("C41").Value = [double]"4.303008943067832e-32"
$ws.Range("C42").Value = [double]"1.726718737528276e-28"
$ws.Range("C43").Value = [double]"6.119578015645701e-23"
$ws.Range("C44").Value = [double]"6.550101718947557e-19"
$ws.Range("C45").Value = [double]"7.336158574362878e-16"
$ws.Range("C46").Value = [double]"1.804403142916507e-13"
$ws.Range("C47").Value = [double]"1.530379876013285e-11"
$ws.Range("C48").Value = [double]"5.004985242095217e-08"
$ws.Range("C49").Value = [double]"1.183890786452967e-05"
$ws.Range("C50").Value = [double]"0.0006019780874969857"
$ws.Range("C51").Value = [double]"0.01152179769393968"
$ws.Range("C52").Value = [double]"0.7121304377251496"
$ws.Range("C53").Value = [double]"10.92115242256195"
$ws.Range("C54").Value = [double]"75.33112781639041"
$ws.Range("C55").Value = [double]"315.3671769034067"
$ws.Range("C56").Value = [double]"2254.663230382765"
$ws.Range("C57").Value = [double]"8054.000456659141"
$ws.Range("C58").Value = [double]"19420.0782703348"
$ws.Range("C59").Value = [double]"36723.722966535"
$ws.Range("C60").Value = [double]"59174.89451091763"
$ws.Range("C61").Value = [double]"85374.39430066453"
